$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "45766004"
$ws.Range("D16").Value = "GLORIA BEATRIZ ARRIETA MEDINA"
$ws.Range("E16").Value = "1908"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 877803
$ws.Range("C17").Value = "45766004"
$ws.Range("D17").Value = "GLORIA BEATRIZ ARRIETA MEDINA"
$ws.Range("E17").Value = "1907"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 877803
$ws.Range("C18").Value = "45766004"
$ws.Range("D18").Value = "GLORIA BEATRIZ ARRIETA MEDINA"
$ws.Range("E18").Value = "1906"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 877803
$ws.Range("C19").Value = "20204308"
$ws.Range("D19").Value = "HEYDY CAROLINA ORTIZ CABRERA"
$ws.Range("E19").Value = "1905"
$ws.Range("F19").Value = 36000
$ws.Range("G19").Value = 900000
$ws.Range("C20").Value = "1030646734"
$ws.Range("D20").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E20").Value = "2105"
$ws.Range("F20").Value = 28800
$ws.Range("G20").Value = 900000
$ws.Range("C21").Value = "1030646734"
$ws.Range("D21").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E21").Value = "2104"
$ws.Range("F21").Value = 36000
$ws.Range("G21").Value = 900000
$ws.Range("C22").Value = "1030646734"
$ws.Range("D22").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E22").Value = "2103"
$ws.Range("F22").Value = 36000
$ws.Range("G22").Value = 900000
$ws.Range("C23").Value = "1030646734"
$ws.Range("D23").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E23").Value = "2102"
$ws.Range("F23").Value = 36000
$ws.Range("G23").Value = 900000
$ws.Range("C24").Value = "1030646734"
$ws.Range("D24").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E24").Value = "2101"
$ws.Range("F24").Value = 36000
$ws.Range("G24").Value = 900000
$ws.Range("C25").Value = "1030646734"
$ws.Range("D25").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E25").Value = "2012"
$ws.Range("F25").Value = 36000
$ws.Range("G25").Value = 900000
$ws.Range("C26").Value = "1030646734"
$ws.Range("D26").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E26").Value = "2011"
$ws.Range("F26").Value = 36000
$ws.Range("G26").Value = 900000
$ws.Range("C27").Value = "1030646734"
$ws.Range("D27").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E27").Value = "2010"
$ws.Range("F27").Value = 36000
$ws.Range("G27").Value = 900000
$ws.Range("C28").Value = "1030646734"
$ws.Range("D28").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E28").Value = "2009"
$ws.Range("F28").Value = 36000
$ws.Range("G28").Value = 900000
$ws.Range("C29").Value = "1030646734"
$ws.Range("D29").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E29").Value = "2008"
$ws.Range("F29").Value = 36000
$ws.Range("G29").Value = 900000
$ws.Range("C30").Value = "1030646734"
$ws.Range("D30").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E30").Value = "2007"
$ws.Range("F30").Value = 36000
$ws.Range("G30").Value = 900000
$ws.Range("C31").Value = "1030646734"
$ws.Range("D31").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E31").Value = "2006"
$ws.Range("F31").Value = 36000
$ws.Range("G31").Value = 900000
$ws.Range("C32").Value = "1030646734"
$ws.Range("D32").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E32").Value = "2005"
$ws.Range("F32").Value = 36000
$ws.Range("G32").Value = 900000
$ws.Range("C33").Value = "1030646734"
$ws.Range("D33").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E33").Value = "2004"
$ws.Range("F33").Value = 36000
$ws.Range("G33").Value = 900000
$ws.Range("C34").Value = "1030646734"
$ws.Range("D34").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E34").Value = "2003"
$ws.Range("F34").Value = 36000
$ws.Range("G34").Value = 900000
$ws.Range("C35").Value = "1030646734"
$ws.Range("D35").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E35").Value = "2002"
$ws.Range("F35").Value = 36000
$ws.Range("G35").Value = 900000
$ws.Range("C36").Value = "1030646734"
$ws.Range("D36").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E36").Value = "2001"
$ws.Range("F36").Value = 36000
$ws.Range("G36").Value = 900000
$ws.Range("C37").Value = "1030646734"
$ws.Range("D37").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E37").Value = "1912"
$ws.Range("F37").Value = 36000
$ws.Range("G37").Value = 900000
$ws.Range("C38").Value = "1030646734"
$ws.Range("D38").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E38").Value = "1911"
$ws.Range("F38").Value = 36000
$ws.Range("G38").Value = 900000
$ws.Range("C39").Value = "1030646734"
$ws.Range("D39").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E39").Value = "1910"
$ws.Range("F39").Value = 36000
$ws.Range("G39").Value = 900000
$ws.Range("C40").Value = "1030646734"
$ws.Range("D40").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E40").Value = "1909"
$ws.Range("F40").Value = 36000
$ws.Range("G40").Value = 900000
$ws.Range("C41").Value = "1030646734"
$ws.Range("D41").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E41").Value = "1908"
$ws.Range("F41").Value = 36000
$ws.Range("G41").Value = 900000
$ws.Range("C42").Value = "1030646734"
$ws.Range("D42").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E42").Value = "1907"
$ws.Range("F42").Value = 36000
$ws.Range("G42").Value = 900000
$ws.Range("C43").Value = "1030646734"
$ws.Range("D43").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E43").Value = "1906"
$ws.Range("F43").Value = 36000
$ws.Range("G43").Value = 900000
$ws.Range("C44").Value = "1030646734"
$ws.Range("D44").Value = "PAOLA ANDREA ARENAS MOGOLLON"
$ws.Range("E44").Value = "1905"
$ws.Range("F44").Value = 36000
$ws.Range("G44").Value = 900000